# Fruta / hortaliza, semanal
# Weekly refresh: prepend the newest week's two quality-grade rows
# ("Pintón" / "Primera Pintón") at row 522, pushing every existing row
# down by two. The two oldest rows that fall off the end of the former
# range (old rows 616/617) land as the new rows 618/619.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at row 522 - this shifts rows 522:617 down
# to 524:619 (values + formatting carried along automatically).
$ws.Rows.Item(522).Insert()
$ws.Rows.Item(522).Insert()

# Columns A,B,C,E,F,G,H,I,J,K,Q,R,T are constant for every data row in
# this sheet (same market / product template), so reuse them verbatim
# for the two freshly-inserted rows.
$marketId   = 7
$market     = "Terminal Hortofrutícola Agro Chillán"
$region     = "Ñuble"
$codreg     = 16
$tipo       = "Fruta"
$prodId     = 100108
$producto   = "Tropicales y subtropicales"
$catId      = 100108006
$categoria  = "Plátano"
$variedad   = "Sin especificar"
$unidad     = "`$/caja 20 kilos"
$origen     = "Ecuador"
$kgUnidad   = 20

function Set-DataRow {
    param(
        [int]$row,
        [double]$fecha,
        [string]$calidad,
        [double]$volumen,
        [double]$precioMin,
        [double]$precioMax,
        [double]$precioProm,
        [double]$precioKg
    )

    $ws.Cells.Item($row, 1).Value  = $marketId
    $ws.Cells.Item($row, 2).Value  = $market
    $ws.Cells.Item($row, 3).Value  = $region
    $ws.Cells.Item($row, 4).Value  = $fecha
    $ws.Cells.Item($row, 5).Value  = $codreg
    $ws.Cells.Item($row, 6).Value  = $tipo
    $ws.Cells.Item($row, 7).Value  = $prodId
    $ws.Cells.Item($row, 8).Value  = $producto
    $ws.Cells.Item($row, 9).Value  = $catId
    $ws.Cells.Item($row, 10).Value = $categoria
    $ws.Cells.Item($row, 11).Value = $variedad
    $ws.Cells.Item($row, 12).Value = $calidad
    $ws.Cells.Item($row, 13).Value = $volumen
    $ws.Cells.Item($row, 14).Value = $precioMin
    $ws.Cells.Item($row, 15).Value = $precioMax
    $ws.Cells.Item($row, 16).Value = $precioProm
    $ws.Cells.Item($row, 17).Value = $unidad
    $ws.Cells.Item($row, 18).Value = $origen
    $ws.Cells.Item($row, 19).Value = $precioKg
    $ws.Cells.Item($row, 20).Value = $kgUnidad
}

Set-DataRow -row 522 -fecha 44694 -calidad "Pintón"          -volumen 80  -precioMin 11000 -precioMax 11000 -precioProm 11000 -precioKg 550
Set-DataRow -row 523 -fecha 44694 -calidad "Primera Pintón"  -volumen 160 -precioMin 12000 -precioMax 13000 -precioProm 12500 -precioKg 625

Write-Output "rows inserted and populated"
